$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '47.596.65'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.00%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.489.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.67%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '105.21'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.526'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.10%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.543'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.43'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0815'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.56%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.124'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.30'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.17'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.875.24'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.493.98'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.845'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '47.456.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.89%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0937'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '70.75'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '251.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.83%  '
$ws.Range("E24").Value = '  +5.88%  '
$ws.Range("E25").Value = '  +3.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.15'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.48%  '
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.87%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.21'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.05'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.78%  '
$ws.Range("E31").Value = '  +7.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.95%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.83'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.71%  '
$ws.Range("E34").Value = '  +2.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0783'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.12%  '
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.96'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.62'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.22%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.99'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.44%  '
$ws.Range("E40").Value = '  +1.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.112'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.93%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '121.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.06'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0298'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.81%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.965.70'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.93%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.98'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.22'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.80'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.31'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +12.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.52'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.86%  '
